$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value2 = "Do you currently live in the United States?"
$ws.Range("B20").Value2 = "Yes"
$ws.Range("A20").Font.Bold = $true
$ws.Range("A20").Font.Size = 12
$ws.Range("A20").Font.Color = 0
$ws.Range("B20").Font.Size = 11
$ws.Range("B20").Font.Color = 0
$ws.Range("B20").HorizontalAlignment = -4152
Write-Host "done"
